$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "E-Mail" column (old column E) to make room
# for the new "Mật khẩu (*)" header; this shifts E-Mail -> F and
# Số điện thoại -> G.
$ws.Columns.Item(5).Insert()

# Give the newly inserted column the same width as column D (both end up
# 15.88671875 wide).
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Update header labels to the new "(*)" required-field wording and add the
# new password column header.
$ws.Range("B1").Value = "Tên đăng nhập (*)"
$ws.Range("C1").Value = "Họ (*)"
$ws.Range("D1").Value = "Tên (*)"
$ws.Range("E1").Value = "Mật khẩu (*)"

# Style the whole header row: bold text with a thin box border around every
# cell. Build the combined format on an out-of-the-way helper cell first and
# stamp it onto the header range in one shot via PasteSpecial so the header
# cells land directly on the final (bold + bordered) style.
$helper = $ws.Range("Z1")
$helper.Font.Bold = $true
$helper.Borders.LineStyle = 1
$helper.Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)
$helper.Clear()

# Move the active selection to C5 (matches the saved workbook view).
[void]$ws.Range("C5").Select()
